$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.855.99'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.40%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.239.41'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.43%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.33%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.16'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.47%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '99.33'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.16%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.571'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.79%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.24%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.540'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -3.62%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.32'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.10%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0820'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.97%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.43'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -3.30%  '

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.90%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.579.55'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.49%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.842'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.49%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.15'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.58%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.236.02'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.36%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.772.81'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.43%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.32'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -6.47%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0968'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.12%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.38'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -3.08%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '64.97'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.28%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.07'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -4.52%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '232.88'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.72%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.05'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -5.14%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.19%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.45'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +3.92%  '

$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.19'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.32%  '

$ws.Range("B29").Value = 'InjectiveProtocol'
$ws.Range("C29").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '37.77'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +3.70%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.98'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -6.86%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '158.87'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.91%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '19.95'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.48%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0835'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.70%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.67'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.64%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.14'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.11%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +6.35%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.91'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.22%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.117'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -2.67%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '15.89'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +11.75%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.63'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.82%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.12'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -6.46%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0310'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.83%  '

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.28%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.758.82'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.02%  '

$ws.Range("B45").Value = 'ordi'
$ws.Range("C45").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '73.98'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.37%  '

$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.193'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -4.60%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '80.54'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -4.22%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.10'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -3.46%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '102.41'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.71%  '

$ws.Range("B50").Value = 'MultiversX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '56.91'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.12%  '

$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.63'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.11%  '
